# Update data package: shorten the "AUXÍLIOS ..." description text and
# move the active selection from D17 to K3 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell whose shared-string text changed from
# "RECURSOS RECEBIDOS PARA AUXÍLIOS DOENÇA, FUNERAL, ALIMENTAÇÃO, TRANSPORTE E FARDAMENTO"
# to "RECURSOS RECEBIDOS PARA AUXÍLIOS".
$ws.Range("B9").Value = "RECURSOS RECEBIDOS PARA AUXÍLIOS"

# Move the sheet's selection/active cell to K3 (was D17).
$ws.Range("K3").Select()
